$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.652.59"
$ws.Range('E2').Value2 = '  -0.62%  '
$ws.Range('D3').Value = "'2.541.69"
$ws.Range('E3').Value2 = '  -0.20%  '
$ws.Range('E4').Value2 = '  +0.02%  '
$ws.Range('D5').Value = "'314.13"
$ws.Range('E5').Value2 = '  +3.10%  '
$ws.Range('D6').Value = "'95.77"
$ws.Range('E6').Value2 = '  -2.44%  '
$ws.Range('E7').Value2 = '  +0.25%  '
$ws.Range('E8').Value2 = '  -0.04%  '
$ws.Range('E9').Value2 = '  -1.68%  '
$ws.Range('D10').Value = "'36.53"
$ws.Range('E10').Value2 = '  -1.05%  '
$ws.Range('E11').Value2 = '  -2.17%  '
$ws.Range('E12').Value2 = '  -0.29%  '
$ws.Range('E13').Value2 = '  -1.18%  '
$ws.Range('D14').Value = "'2.932.38"
$ws.Range('E14').Value2 = '  -0.24%  '
$ws.Range('D15').Value = "'15.78"
$ws.Range('E15').Value2 = '  +4.35%  '
$ws.Range('D16').Value = "'2.539.44"
$ws.Range('E16').Value2 = '  +1.23%  '
$ws.Range('D17').Value = "'0.870"
$ws.Range('E17').Value2 = '  -0.67%  '
$ws.Range('D18').Value = "'42.717.96"
$ws.Range('E18').Value2 = '  -0.57%  '
$ws.Range('D19').Value = "'13.15"
$ws.Range('E19').Value2 = '  -4.91%  '
$ws.Range('E20').Value2 = '  +1.34%  '
$ws.Range('E21').Value2 = '  -2.55%  '
$ws.Range('D22').Value = "'71.25"
$ws.Range('E22').Value2 = '  -0.93%  '
$ws.Range('D23').Value = "'255.02"
$ws.Range('E23').Value2 = '  -0.32%  '
$ws.Range('D24').Value = "'2.96"
$ws.Range('E24').Value2 = '  -0.05%  '
$ws.Range('D25').Value = "'2.05"
$ws.Range('E25').Value2 = '  -1.95%  '
$ws.Range('D26').Value = "'27.63"
$ws.Range('E26').Value2 = '  -1.64%  '
$ws.Range('D27').Value = "'1.00"
$ws.Range('E27').Value2 = '  +0.14%  '
$ws.Range('E28').Value2 = '  +11.98%  '
$ws.Range('D29').Value = "'39.76"
$ws.Range('E29').Value2 = '  +5.27%  '
$ws.Range('E30').Value2 = '  -1.13%  '
$ws.Range('D31').Value = "'5.94"
$ws.Range('E31').Value2 = '  -4.15%  '
$ws.Range('D32').Value = "'155.88"
$ws.Range('E32').Value2 = '  -1.45%  '
$ws.Range('D33').Value = "'19.88"
$ws.Range('E33').Value2 = '  +0.61%  '
$ws.Range('D34').Value = "'3.41"
$ws.Range('E34').Value2 = '  +2.66%  '
$ws.Range('E35').Value2 = '  +0.08%  '
$ws.Range('E36').Value2 = '  -1.27%  '
$ws.Range('E37').Value2 = '  -0.13%  '
$ws.Range('B38').Value2 = 'EnergySwap'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = "'25.57"
$ws.Range('E38').Value2 = '  +0.94%  '
$ws.Range('B39').Value2 = 'Kaspa'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = "'0.113"
$ws.Range('E39').Value2 = '  -3.43%  '
$ws.Range('E40').Value2 = '  -0.23%  '
$ws.Range('D41').Value = "'2.29"
$ws.Range('E41').Value2 = '  +10.63%  '
$ws.Range('E42').Value2 = '  -0.89%  '
$ws.Range('D43').Value = "'3.86"
$ws.Range('E43').Value2 = '  -1.36%  '
$ws.Range('E44').Value2 = '  -0.69%  '
$ws.Range('E45').Value2 = '  +0.01%  '
$ws.Range('D46').Value = "'2.048.89"
$ws.Range('E46').Value2 = '  -2.69%  '
$ws.Range('D47').Value = "'85.87"
$ws.Range('E47').Value2 = '  -1.22%  '
$ws.Range('D48').Value = "'8.93"
$ws.Range('E48').Value2 = '  -0.80%  '
$ws.Range('B49').Value2 = 'RocketPoolETH'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = "'2.782.28"
$ws.Range('E49').Value2 = '  -0.55%  '
$ws.Range('B50').Value2 = 'ordi'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').Value = "'74.65"
$ws.Range('E50').Value2 = '  -0.33%  '
$ws.Range('E51').Value2 = '  -0.42%  '
